$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = "267.31"
    3  = "22.87"
    4  = "6.358"
    5  = "0.06194"
    6  = "3.603"
    7  = "6.687"
    8  = "1.396"
    9  = "0.8307"
    10 = "0.01362"
    11 = "0.1611"
    12 = "0.08224"
    13 = "0.03400"
    14 = "0.03158"
    15 = "0.09289"
    17 = "0.001709"
    18 = "0.04870"
    20 = "0.005380"
    21 = "0.001095"
    22 = "0.0001505"
    23 = "3.762"
    24 = "2.368"
    25 = "0.3344"
    26 = "0.1247"
    27 = "0.0002696"
    40 = "0.04679"
    41 = "0.006932"
    42 = "0.1155"
    43 = "0.003342"
    44 = "0.01237"
    45 = "0.00006279"
    46 = "0.00000000754"
    47 = "0.7930"
    48 = "0.1555"
    49 = "0.00002111"
    50 = "0.01247"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
}
